$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2307
    $ws.Range("F3").Value = 1760
    $ws.Range("F4").Value = 342
    $ws.Range("F5").Value = 1101
    $ws.Range("F6").Value = 942
    $ws.Range("F8").Value = 5871
}
